$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap data between row 208 and row 211
$ws.Range("B208").Value = 7404212
$ws.Range("B211").Value = 7404216
$ws.Range("F208").Value = "Envigado FC"
$ws.Range("F211").Value = "Independiente Santa Fe"
$ws.Range("G208").Value = "Deportivo Pasto"
$ws.Range("G211").Value = "Once Caldas"
$ws.Range("H208").Value = 1
$ws.Range("H211").Value = 0
$ws.Range("J208").Value = "D"
$ws.Range("J211").Value = "A"
$ws.Range("K208").Value = 2.6
$ws.Range("K211").Value = 1.85
$ws.Range("L208").Value = 2.875
$ws.Range("L211").Value = 3.1
$ws.Range("M208").Value = 2.8
$ws.Range("M211").Value = 4.2
$ws.Range("N208").Value = 2.8
$ws.Range("N211").Value = 2.25
$ws.Range("P208").Value = 2.625
$ws.Range("P211").Value = 3.3
$ws.Range("Q208").Value = 0
$ws.Range("Q211").Value = -0.25
$ws.Range("R208").Value = 1.975
$ws.Range("R211").Value = 1.9
$ws.Range("S208").Value = 1.875
$ws.Range("S211").Value = 1.9
$ws.Range("U208").Value = 2.025
$ws.Range("U211").Value = 1.925
$ws.Range("V208").Value = 1.825
$ws.Range("V211").Value = 1.925
$ws.Range("X208").Value = 2.2
$ws.Range("X211").Value = -1
$ws.Range("Y208").Value = -1
$ws.Range("Y211").Value = 2.3
$ws.Range("Z208").Value = 0
$ws.Range("Z211").Value = -1
$ws.Range("AA208").Value = -0
$ws.Range("AA211").Value = 0.8999999999999999
$ws.Range("AC208").Value = 0.825
$ws.Range("AC211").Value = 0.925

# Swap data between row 209 and row 212
$ws.Range("B209").Value = 7404214
$ws.Range("B212").Value = 7404217
$ws.Range("F209").Value = "Boyaca Chico"
$ws.Range("F212").Value = "Alianza Petrolera"
$ws.Range("G209").Value = "Deportivo Cali"
$ws.Range("G212").Value = "Deportivo Pereira"
$ws.Range("H209").Value = 1
$ws.Range("H212").Value = 2
$ws.Range("J209").Value = "D"
$ws.Range("J212").Value = "H"
$ws.Range("K209").Value = 3.2
$ws.Range("K212").Value = 1.95
$ws.Range("L209").Value = 3.1
$ws.Range("L212").Value = 3.2
$ws.Range("M209").Value = 2.2
$ws.Range("M212").Value = 3.75
$ws.Range("N209").Value = 3.6
$ws.Range("N212").Value = 1.95
$ws.Range("O209").Value = 3
$ws.Range("O212").Value = 3.2
$ws.Range("P209").Value = 2.25
$ws.Range("P212").Value = 4.75
$ws.Range("Q209").Value = 0.25
$ws.Range("Q212").Value = -0.5
$ws.Range("R209").Value = 1.95
$ws.Range("R212").Value = 1.925
$ws.Range("S209").Value = 1.9
$ws.Range("S212").Value = 1.875
$ws.Range("T209").Value = 2.25
$ws.Range("T212").Value = 2
$ws.Range("U209").Value = 1.875
$ws.Range("U212").Value = 1.825
$ws.Range("W209").Value = -1
$ws.Range("W212").Value = 0.95
$ws.Range("X209").Value = 2
$ws.Range("X212").Value = -1
$ws.Range("Z209").Value = 0.475
$ws.Range("Z212").Value = 0.925
$ws.Range("AA209").Value = -0.5
$ws.Range("AA212").Value = -1
$ws.Range("AB209").Value = -0.5
$ws.Range("AB212").Value = 0.825
$ws.Range("AC209").Value = 0.4875
$ws.Range("AC212").Value = -1

# Swap data between row 240 and row 241
$ws.Range("B240").Value = 7528603
$ws.Range("B241").Value = 7528135
$ws.Range("F240").Value = "Junior"
$ws.Range("F241").Value = "Independiente Medellin"
$ws.Range("G240").Value = "Deportes Tolima"
$ws.Range("G241").Value = "America de Cali"
$ws.Range("H240").Value = 4
$ws.Range("H241").Value = 2
$ws.Range("I240").Value = 2
$ws.Range("I241").Value = 1
$ws.Range("K240").Value = 1.95
$ws.Range("K241").Value = 2.15
$ws.Range("M240").Value = 4
$ws.Range("M241").Value = 3.4
$ws.Range("N240").Value = 1.909
$ws.Range("N241").Value = 2.375
$ws.Range("O240").Value = 3.75
$ws.Range("O241").Value = 3.3
$ws.Range("P240").Value = 3.8
$ws.Range("P241").Value = 3.1
$ws.Range("Q240").Value = -0.5
$ws.Range("Q241").Value = -0.25
$ws.Range("R240").Value = 1.9
$ws.Range("R241").Value = 2
$ws.Range("S240").Value = 1.9
$ws.Range("S241").Value = 1.8
$ws.Range("U240").Value = 1.85
$ws.Range("U241").Value = 1.975
$ws.Range("V240").Value = 1.95
$ws.Range("V241").Value = 1.825
$ws.Range("W240").Value = 0.909
$ws.Range("W241").Value = 1.375
$ws.Range("Z240").Value = 0.8999999999999999
$ws.Range("Z241").Value = 1
$ws.Range("AB240").Value = 0.8500000000000001
$ws.Range("AB241").Value = 0.9750000000000001

# Row 373 single-cell edits
$ws.Range("N373").Value = 2.875
$ws.Range("P373").Value = 2.45
$ws.Range("Q373").Value = 0
$ws.Range("R373").Value = 2.1
$ws.Range("S373").Value = 1.775

# Row 374 single-cell edits
$ws.Range("R374").Value = 1.85
$ws.Range("S374").Value = 2

# Row 379 single-cell edits
$ws.Range("R379").Value = 1.925
$ws.Range("S379").Value = 1.925
